$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2404091954231262
$ws.Range("B1").Value = 0.2316491901874542
$ws.Range("C1").Value = 0.241667851805687
$ws.Range("D1").Value = 0.3276466131210327
$ws.Range("E1").Value = 0.5902460217475891
